$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# res_bus/vm_pu.xlsx update: case re-run with 380 kV slack bus (1.02 pu)
# Each data row (2-25) gets new per-unit voltage magnitudes for B:F and I:N;
# column G (slack, always 1) and the empty column H are left untouched.

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.004727430580851
$bf[0,2] = 1.030263256528299
$bf[0,3] = 1.007762213487518
$bf[0,4] = 1.002867825820396
$ws.Range("B2:F2").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.02955528997503
$in[0,1] = 1.010014891386621
$in[0,2] = 1.033074774661431
$in[0,3] = 1.010640516912302
$in[0,4] = 1.005761089072324
$in[0,5] = 1.007465362557762
$ws.Range("I2:N2").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.005822002272901
$bf[0,2] = 1.030606612076188
$bf[0,3] = 1.008694518224989
$bf[0,4] = 1.004617363073518
$ws.Range("B3:F3").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029546965970127
$in[0,1] = 1.010738946473879
$in[0,2] = 1.033227698754684
$in[0,3] = 1.01137594559515
$in[0,4] = 1.007310333110168
$in[0,5] = 1.007712251451404
$ws.Range("I3:N3").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.006530092566826
$bf[0,2] = 1.030828240658591
$bf[0,3] = 1.009297975763954
$bf[0,4] = 1.005748917304874
$ws.Range("B4:F4").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029539753632458
$in[0,1] = 1.01120681070201
$in[0,2] = 1.0333253549493
$in[0,3] = 1.011851397160795
$in[0,4] = 1.008311870664742
$in[0,5] = 1.007871580600961
$ws.Range("I4:N4").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.006827735276856
$bf[0,2] = 1.030921280219724
$bf[0,3] = 1.009551717083347
$bf[0,4] = 1.006224509550559
$ws.Range("B5:F5").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029536283262584
$in[0,1] = 1.011403347421962
$in[0,2] = 1.033366097555123
$in[0,3] = 1.012051177842443
$in[0,4] = 1.008732703681971
$in[0,5] = 1.007938461346076
$ws.Range("I5:N5").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.006877708613519
$bf[0,2] = 1.030936894093018
$bf[0,3] = 1.00959432420678
$bf[0,4] = 1.006304357228488
$ws.Range("B6:F6").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029535674837201
$in[0,1] = 1.011436337855459
$in[0,2] = 1.033372920065691
$in[0,3] = 1.012084716123469
$in[0,4] = 1.008803351132689
$in[0,5] = 1.007949684992606
$ws.Range("I6:N6").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.006534069831341
$bf[0,2] = 1.030829484383945
$bf[0,3] = 1.009301366078558
$bf[0,4] = 1.005755272624287
$ws.Range("B7:F7").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029539708984861
$in[0,1] = 1.01120943743596
$in[0,2] = 1.033325900582478
$in[0,3] = 1.011854067027955
$in[0,4] = 1.008317494683953
$in[0,5] = 1.007872474662037
$ws.Range("I7:N7").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.005097383039582
$bf[0,2] = 1.030379406564799
$bf[0,3] = 1.008077250859506
$bf[0,4] = 1.003459202297312
$ws.Range("B8:F8").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.02955285428455
$in[0,1] = 1.010259724350782
$in[0,2] = 1.033126723317437
$in[0,3] = 1.010889146272544
$in[0,4] = 1.006284860424486
$in[0,5] = 1.00754888794167
$ws.Range("I8:N8").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.002564339795591
$bf[0,2] = 1.029582245551909
$bf[0,3] = 1.005921629741186
$bf[0,4] = 0.9994088644037524
$ws.Range("B9:F9").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029562081003761
$in[0,1] = 1.008581170752679
$in[0,2] = 1.032765901988229
$in[0,3] = 1.009185553564632
$in[0,4] = 1.002695617146474
$in[0,5] = 1.00697541805312
$ws.Range("I9:N9").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.000874535888566
$bf[0,2] = 1.029048220613477
$bf[0,3] = 1.004485421829613
$bf[0,4] = 0.9967050907846262
$ws.Range("B10:F10").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029558927936093
$in[0,1] = 1.007458639509106
$in[0,2] = 1.032518842347417
$in[0,3] = 1.008047530887223
$in[0,4] = 1.000297222101511
$in[0,5] = 1.006590881253872
$ws.Range("I10:N10").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.000142534010965
$bf[0,2] = 1.028816400681362
$bf[0,3] = 1.003863715011103
$bf[0,4] = 0.995533349279506
$ws.Range("B11:F11").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029555369739793
$in[0,1] = 1.006971717307854
$in[0,2] = 1.032410340936442
$in[0,3] = 1.007554189785336
$in[0,4] = 0.999257251006217
$in[0,5] = 1.006423839134035
$ws.Range("I11:N11").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 0.99987058688519
$bf[0,2] = 1.028730207336787
$bf[0,3] = 1.003632810757945
$bf[0,4] = 0.9950979518761309
$ws.Range("B12:F12").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.02955371971951
$in[0,1] = 1.006790721502815
$in[0,2] = 1.032369811708957
$in[0,3] = 1.007370853485298
$in[0,4] = 0.9988707312407346
$in[0,5] = 1.006361711295728
$ws.Range("I12:N12").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 0.9999289227360615
$bf[0,2] = 1.028748699911505
$bf[0,3] = 1.00368233938652
$bf[0,4] = 0.9951913535466643
$ws.Range("B13:F13").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029554088495152
$in[0,1] = 1.006829551690842
$in[0,2] = 1.032378515607504
$in[0,3] = 1.007410183731069
$in[0,4] = 0.9989536515352102
$in[0,5] = 1.006375041600458
$ws.Range("I13:N13").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.000120055795136
$bf[0,2] = 1.028809277633772
$bf[0,3] = 1.003844627887068
$bf[0,4] = 0.9954973625354055
$ws.Range("B14:F14").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029555240040407
$in[0,1] = 1.006956758827764
$in[0,2] = 1.032406995397466
$in[0,3] = 1.007539036944252
$in[0,4] = 0.9992253058768764
$in[0,5] = 1.006418705282296
$ws.Range("I14:N14").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.000237812589679
$bf[0,2] = 1.028846590368057
$bf[0,3] = 1.003944622492043
$bf[0,4] = 0.9956858831607941
$ws.Range("B15:F15").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029555906069241
$in[0,1] = 1.007035117889274
$in[0,2] = 1.032424512714591
$in[0,3] = 1.007618416018253
$in[0,4] = 0.9993926505345242
$in[0,5] = 1.006445597173323
$ws.Range("I15:N15").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.000923109737314
$bf[0,2] = 1.029063593663635
$bf[0,3] = 1.004526686116763
$bf[0,4] = 0.9967828335567679
$ws.Range("B16:F16").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029559117975007
$in[0,1] = 1.007490936651406
$in[0,2] = 1.032526011298701
$in[0,3] = 1.008080260149639
$in[0,4] = 1.000366210124534
$in[0,5] = 1.006601955966833
$ws.Range("I16:N16").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.001352894509748
$bf[0,2] = 1.029199559676678
$bf[0,3] = 1.00489184635577
$bf[0,4] = 0.997470648322167
$ws.Range("B17:F17").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029560546404202
$in[0,1] = 1.007776628254542
$in[0,2] = 1.032589272239751
$in[0,3] = 1.00836980892789
$in[0,4] = 1.000976502807035
$in[0,5] = 1.006699892041016
$ws.Range("I17:N17").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.00160355141928
$bf[0,2] = 1.029278809869972
$bf[0,3] = 1.005104855702431
$bf[0,4] = 0.9978717443451044
$ws.Range("B18:F18").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029561167926103
$in[0,1] = 1.007943184624868
$in[0,2] = 1.032626024211051
$in[0,3] = 1.008538642982987
$in[0,4] = 1.001332336979295
$in[0,5] = 1.006756964900417
$ws.Range("I18:N18").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.001689014075675
$bf[0,2] = 1.029305822478735
$bf[0,3] = 1.005177489416966
$bf[0,4] = 0.9980084920162302
$ws.Range("B19:F19").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029561343914426
$in[0,1] = 1.007999962080195
$in[0,2] = 1.032638530682611
$in[0,3] = 1.008596201784786
$in[0,4] = 1.001453643932613
$in[0,5] = 1.006776416535688
$ws.Range("I19:N19").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.001306785723031
$bf[0,2] = 1.029184977633886
$bf[0,3] = 1.004852666315862
$bf[0,4] = 0.9973968621890754
$ws.Range("B20:F20").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029560415031944
$in[0,1] = 1.007745984820907
$in[0,2] = 1.032582500137955
$in[0,3] = 1.008338748771612
$in[0,4] = 1.000911038655298
$in[0,5] = 1.006689389768212
$ws.Range("I20:N20").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.000063773251871
$bf[0,2] = 1.02879144133297
$bf[0,3] = 1.003796837283646
$bf[0,4] = 0.9954072550388525
$ws.Range("B21:F21").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029554909994162
$in[0,1] = 1.006919303136536
$in[0,2] = 1.032398615060048
$in[0,3] = 1.007501095313476
$in[0,4] = 0.9991453168265775
$in[0,5] = 1.006405849655614
$ws.Range("I21:N21").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 0.9992819567113377
$bf[0,2] = 1.0285435183423
$bf[0,3] = 1.003133141372291
$bf[0,4] = 0.9941553759467042
$ws.Range("B22:F22").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029549549582982
$in[0,1] = 1.006398775255053
$in[0,2] = 1.032281686918611
$in[0,3] = 1.006973922460185
$in[0,4] = 0.9980338120025594
$in[0,5] = 1.006227108057501
$ws.Range("I22:N22").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 0.9996964405381674
$bf[0,2] = 1.028674992704325
$bf[0,3] = 1.003484965928803
$bf[0,4] = 0.9948191132268769
$ws.Range("B23:F23").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029552570889779
$in[0,1] = 1.006674789740152
$in[0,2] = 1.032343796520245
$in[0,3] = 1.007253435431472
$in[0,4] = 0.9986231708395156
$in[0,5] = 1.006321906973569
$ws.Range("I23:N23").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.001327620389373
$bf[0,2] = 1.029191566806252
$bf[0,3] = 1.004870370035164
$bf[0,4] = 0.9974302032580508
$ws.Range("B24:F24").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029560475047558
$in[0,1] = 1.007759831525004
$in[0,2] = 1.032585560613934
$in[0,3] = 1.008352783688695
$in[0,4] = 1.000940619516847
$in[0,5] = 1.006694135452648
$ws.Range("I24:N24").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.003219376118636
$bf[0,2] = 1.029788797171083
$bf[0,3] = 1.006478749260445
$bf[0,4] = 1.000456554253951
$ws.Range("B25:F25").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0] = 1.029561340424493
$in[0,1] = 1.009015725207216
$in[0,2] = 1.032860338570589
$in[0,3] = 1.009626370872034
$in[0,4] = 1.0036244642001
$in[0,5] = 1.007124063766357
$ws.Range("I25:N25").Value = $in

